# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.768.37'
$ws.Range('E2').Value = '  -2.78%  '
$ws.Range('D3').Value = '2.091.33'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '345.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.008'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4482'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09364'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.03'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.178'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.11'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.92%  '
$ws.Range('D13').Value = '2.101.50'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.770'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.088'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '99.12'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001161'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.010'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '20.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06687'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.008'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.186'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.27%  '
$ws.Range('D23').Value = '29.849.20'
$ws.Range('E23').Value = '  -2.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.70'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.315'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.04%  '
$ws.Range('D26').Value = '2.344.96'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.542'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('E30').Value = '  -1.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.160'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1055'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.615'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.211'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.951'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.147'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02570'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06748'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2274'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.43'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6878'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.293'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6618'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.25'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.298'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.639'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000344'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.76%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.219'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '81.84'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07131'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.34%  '
